$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "sku"
$ws.Range("B2").Value = "name"
$ws.Range("C2").Value = "quantity"
$ws.Range("D2").Value = "cost_per"
$ws.Range("E2").Value = "total_cost"
